$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# 1) Reposition/resize the rotated round-rect "Rectangle 65" (id=117)
#    old: off x=6213739 y=4560376  ext cx=1371599 cy=328045  (EMU)
#    new: off x=6055893 y=4718220  ext cx=1687291 cy=328045  (EMU)
$rect = $s.Shapes.Item(29)

$rect.Left  = 6055893 / 12700
$rect.Top   = 4718220 / 12700
$rect.Width = 1687291 / 12700

# Height needs an extra nudge: going EMU -> points -> (internal float) -> EMU
# can truncate the last unit, so creep the points value up until the round
# trip lands back on the exact target EMU.
$targetEmu = 328045
$pts = $targetEmu / 12700.0
$rect.Height = $pts
$tries = 0
while ($tries -lt 20) {
    $achievedEmu = [math]::Round($rect.Height * 12700)
    if ($achievedEmu -ge $targetEmu) { break }
    $pts += 0.000001
    $rect.Height = $pts
    $tries++
}

# 2) Merge the "Task" + "ListPanel" runs into a single "TaskListPanel" run.
#    TextRange.Text already reads back as the concatenation of both runs, so
#    assigning that same string is treated as a no-op; first set an unrelated
#    placeholder (sharing no prefix/suffix with the target) to force a full
#    single-run rewrite, then set the real text.
$taskListPanel = $s.Shapes.Item(41)
$taskListPanel.TextFrame.TextRange.Text = "ZZZZZZZZZZZZZ"
$taskListPanel.TextFrame.TextRange.Text = "TaskListPanel"

# 3) Merge the "Task" + "Card" runs into a single "TaskCard" run (same trick)
$taskCard = $s.Shapes.Item(44)
$taskCard.TextFrame.TextRange.Text = "ZZZZZZZZ"
$taskCard.TextFrame.TextRange.Text = "TaskCard"
